# Fruta / hortaliza, semanal
# Insert a new weekly record at row 345 (Choclo, Dulce o Americano, Primera,
# Región de Arica y Parinacota), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(345).Insert()

$ws.Range("A345").Value = 3
$ws.Range("B345").Value = "Femacal de La Calera"
$ws.Range("C345").Value = "Coquimbo"
$ws.Range("D345").Value = 44461
$ws.Range("E345").Value = 5
$ws.Range("F345").Value = 100112024
$ws.Range("G345").Value = "Choclo"
$ws.Range("H345").Value = "Dulce o Americano"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 35
$ws.Range("K345").Value = 32000
$ws.Range("L345").Value = 32000
$ws.Range("M345").Value = 32000
$ws.Range("N345").Value = "`$/malla 70 unidades"
$ws.Range("O345").Value = "Región de Arica y Parinacota"
$ws.Range("P345").Value = 457
$ws.Range("Q345").Value = 70
$ws.Range("R345").Value = "Hortaliza"
